$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.530.89"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "3.538.87"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "609.02"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +4.76%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "173.18"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.617"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").Value = "3.535.22"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +5.87%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.77"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  -1.13%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "47.32"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "4.111.10"
$ws.Range("E15").Value = "  +1.18%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "627.78"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -6.80%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "8.41"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "70.545.97"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "3.542.39"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  -1.72%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "17.39"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.00"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -10.99%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.886"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.36%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "15.92"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.04%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "96.70"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.81%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "3.85"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.61"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.23%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.21"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.15%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "33.43"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.82%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "8.50"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  -2.05%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.99"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.36%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "568.30"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.62%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "3.61"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  -0.42%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "57.59"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("E40").Value = "  +0.26%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0460"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +5.51%  "
$ws.Range("E42").Value = "  +5.45%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.328"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").Value = "3.345.27"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("D46").Value = "0.0₃0712"
$ws.Range("E46").Value = "  +0.95%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "33.09"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("E49").Value = "  -1.95%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "133.61"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.22%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "5.71"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.88%  "
